$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.457.60"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "1.826.24"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  -0.65%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4572"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07823"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9577"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("D13").Value = "1.835.59"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.827"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.048"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06585"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001019"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("D22").Value = "27.449.22"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.281"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.064.69"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.032"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.277"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09363"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9275"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.575"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.198"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.311"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05953"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02175"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.90%  "

$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.143"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5721"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1818"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.886"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5371"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.893"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06862"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("E51").Value = "  -32.61%  "
